$wb = $excel.ActiveWorkbook

# --- RiskScenarios sheet: insert a new header-ish row at row 2 ---
$ws = $wb.Worksheets.Item("RiskScenarios")
$ws.Rows.Item(2).Insert()

$ws.Range("E2").Value = "n"
$ws.Range("B2").Value = "c"
$ws.Range("C2").Value = "c"
$ws.Range("D2").Value = "c"
$ws.Range("A2").Value = "status"

# Set selection / active cell on this sheet to D13, and make this the active (tabSelected) sheet
$ws.Activate()
$ws.Range("D13").Select()

# --- RiskAssessment sheet used to be the active tab; it no longer is ---
# (switching the active sheet above already clears tabSelected on RiskAssessment)

$wb.Save()
